# Update the two Java stack-trace line references that shifted because of
# the AbstractTemplatesTestSuite refactor that shipped with the 3.0.0
# release (the commit bumping the project from 2.0.2 to 3.0.0).
#
#   AbstractTemplatesTestSuite.java:462 -> :480  (prepareoutputAndGenerate)
#   AbstractTemplatesTestSuite.java:372 -> :389  (generation)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "AbstractTemplatesTestSuite.java:462", $true, $false, $false, $false,
    $false, $true, 1, $false, "AbstractTemplatesTestSuite.java:480", 2)

$d.Content.Find.Execute(
    "AbstractTemplatesTestSuite.java:372", $true, $false, $false, $false,
    $false, $true, 1, $false, "AbstractTemplatesTestSuite.java:389", 2)
